# Apply updated crypto price / 1h-volume values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.004.03"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "3.787.52"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'622.07"
$ws.Range("E5").Value = "  +3.87%  "

$ws.Range("D6").Value = "'177.55"
$ws.Range("E6").Value = "  -3.16%  "

$ws.Range("D7").Value = "3.782.36"
$ws.Range("E7").Value = "  +2.80%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +5.10%  "

$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -3.84%  "

$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").Value = "'41.00"
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "'0.0000263"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("D15").Value = "4.411.65"
$ws.Range("E15").Value = "  +2.62%  "

$ws.Range("D16").Value = "3.781.84"
$ws.Range("E16").Value = "  +2.78%  "

$ws.Range("D17").Value = "70.027.97"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").Value = "'16.84"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "'511.09"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").Value = "'9.52"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("D23").Value = "'0.728"
$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("D24").Value = "'87.80"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  +3.71%  "

$ws.Range("D26").Value = "'13.12"
$ws.Range("E26").Value = "  -3.48%  "

$ws.Range("D27").Value = "'10.99"
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("E28").Value = "  +26.79%  "

$ws.Range("D30").Value = "'2.49"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").Value = "'2.87"
$ws.Range("E31").Value = "  +3.89%  "

$ws.Range("E32").Value = "  -4.29%  "

$ws.Range("D33").Value = "'31.39"
$ws.Range("E33").Value = "  -1.68%  "

$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +4.76%  "

$ws.Range("D37").Value = "'6.21"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("E38").Value = "  +5.64%  "

$ws.Range("D39").Value = "'0.333"
$ws.Range("E39").Value = "  -2.71%  "

$ws.Range("D40").Value = "'2.14"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").Value = "'50.99"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").Value = "'44.98"
$ws.Range("E42").Value = "  -2.61%  "

$ws.Range("D43").Value = "'8.75"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "'417.80"
$ws.Range("E44").Value = "  +4.40%  "

$ws.Range("D45").Value = "'2.83"
$ws.Range("E45").Value = "  +2.62%  "

$ws.Range("D46").Value = "3.034.15"
$ws.Range("E46").Value = "  -4.87%  "

$ws.Range("D47").Value = "'0.0363"
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("E48").Value = "  -2.36%  "

$ws.Range("D49").Value = "'138.73"
$ws.Range("E49").Value = "  +2.12%  "

$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("E51").Value = "  +2.11%  "
